# Remove form_id from remaining forms
#
# The "settings" sheet has columns: form_title | form_id | version | style | namespaces
# This removes the form_id column entirely, shifting version/style/namespaces left
# by one column (B/C/D/E -> B/C/D), and keeps everything else (e.g. the form_title
# column and the auto-updating version formula) intact.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("settings")

# The header-row cell comments describe what belongs in each column. Grab the
# text of the comments that currently sit on the columns which will shift left
# (version/style/namespaces) so they can be re-attached to their new home once
# the form_id column is gone - deleting a column does not relocate comments.
$versionComment    = $ws2.Range("C1").Comment.Text()
$styleComment      = $ws2.Range("D1").Comment.Text()
$namespacesComment = $ws2.Range("E1").Comment.Text()

# Delete the whole form_id column (B) on the settings sheet.
$ws2.Columns.Item(2).Delete()

# Re-home the comments that used to belong to columns C/D/E onto their new
# columns B/C/D (in place, so the original comment/author metadata survives),
# and drop the now out-of-range comment that used to live on column E.
[void]$ws2.Range("B1").Comment.Text($versionComment)
[void]$ws2.Range("C1").Comment.Text($styleComment)
[void]$ws2.Range("D1").Comment.Text($namespacesComment)
[void]$ws2.Range("E1").Comment.Delete()

# Selecting a cell on the settings sheet switches focus to it; restore B1 as
# the selection there (matching where the cursor now naturally lands after the
# form_id column disappears) and then reactivate the survey sheet so the
# workbook's active tab is unchanged.
[void]$ws2.Range("B1").Select()
[void]$ws1.Activate()
